$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.339.46'
$ws.Range('E2').Value = '  +2.76%  '
$ws.Range('D3').Value = '3.697.00'
$ws.Range('E3').Value = '  +8.01%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '582.93'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('D6').Value = '177.60'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('D7').Value = '3.685.48'
$ws.Range('E7').Value = '  +7.88%  '
$ws.Range('D8').Value = '0.614'
$ws.Range('E8').Value = '  +3.82%  '
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').Value = '0.199'
$ws.Range('E10').Value = '  +0.39%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').Value = '6.50'
$ws.Range('E11').Value = '  +19.24%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').Value = '0.608'
$ws.Range('E12').Value = '  +4.29%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = '49.06'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = '0.0000286'
$ws.Range('E14').Value = '  +1.96%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '4.296.63'
$ws.Range('E15').Value = '  +8.22%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value = '678.98'
$ws.Range('E16').Value = '  -1.44%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '8.97'
$ws.Range('E17').Value = '  +3.93%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.707.75'
$ws.Range('E18').Value = '  +8.30%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '71.485.57'
$ws.Range('E19').Value = '  +2.92%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').Value = '0.122'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = '17.93'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '11.55'
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('D23').Value = '0.940'
$ws.Range('E23').Value = '  +4.99%  '
$ws.Range('D24').Value = '17.42'
$ws.Range('E24').Value = '  +2.77%  '
$ws.Range('D25').Value = '102.22'
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('D26').Value = '3.98'
$ws.Range('E26').Value = '  +2.37%  '
$ws.Range('E27').Value = '  +5.51%  '
$ws.Range('D28').Value = '10.31'
$ws.Range('E28').Value = '  +7.59%  '
$ws.Range('D29').Value = '35.11'
$ws.Range('E29').Value = '  +5.08%  '
$ws.Range('D30').Value = '9.13'
$ws.Range('E30').Value = '  +4.45%  '
$ws.Range('D31').Value = '7.33'
$ws.Range('E31').Value = '  +5.03%  '
$ws.Range('D32').Value = '4.06'
$ws.Range('E32').Value = '  +10.08%  '
$ws.Range('D33').Value = '589.51'
$ws.Range('E33').Value = '  +3.30%  '
$ws.Range('D34').Value = '11.18'
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('E35').Value = '  +4.89%  '
$ws.Range('D36').Value = '59.06'
$ws.Range('E36').Value = '  +1.39%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('D38').Value = '3.672.61'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('D39').Value = '0.144'
$ws.Range('E39').Value = '  +4.23%  '
$ws.Range('D40').Value = '0.0₃0765'
$ws.Range('E40').Value = '  +5.08%  '
$ws.Range('D41').Value = '35.21'
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('D42').Value = '3.42'
$ws.Range('E42').Value = '  +4.97%  '
$ws.Range('D43').Value = '2.77'
$ws.Range('E43').Value = '  +4.21%  '
$ws.Range('D44').Value = '0.0457'
$ws.Range('E44').Value = '  +9.82%  '
$ws.Range('D45').Value = '0.347'
$ws.Range('E45').Value = '  +4.70%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '3.37'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').Value = '2.86'
$ws.Range('E47').Value = '  +8.21%  '
$ws.Range('E48').Value = '  +3.58%  '
$ws.Range('D49').Value = '1.43'
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D51').Value = '136.27'
$ws.Range('E51').Value = '  +3.37%  '
